# Implement basket-based elective scheduling with common time slots across all branches
# Updates the Section_A and Section_B timetable sheets so elective/basket
# courses (e.g. CS261/CS262/CS263/CS264/MA261/MA262) line up in the same
# common time slots across both sections.

$wb = $excel.ActiveWorkbook

# ---- Section_A (sheet 1) ----
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("D2").Value = "MA262"
$wsA.Range("E2").Value = "MA262"
$wsA.Range("F2").Value = "Free"

$wsA.Range("D3").Value = "Free"
$wsA.Range("E3").Value = "CS261"
$wsA.Range("F3").Value = "CS264"

$wsA.Range("B5").Value = "CS263"
$wsA.Range("C5").Value = "MA261"
$wsA.Range("D5").Value = "CS262"
$wsA.Range("F5").Value = "CS263"

$wsA.Range("D6").Value = "CS264 (Tutorial)"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "CS262"
$wsA.Range("D7").Value = "CS264"
$wsA.Range("F7").Value = "MA261"

$wsA.Range("F8").Value = "CS263 (Tutorial)"

# ---- Section_B (sheet 2) ----
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("C2").Value = "CS262"
$wsB.Range("D2").Value = "MA262"
$wsB.Range("F2").Value = "CS264"

$wsB.Range("C3").Value = "CS264"
$wsB.Range("D3").Value = "CS264"
$wsB.Range("E3").Value = "CS262"
$wsB.Range("F3").Value = "CS263"

$wsB.Range("C5").Value = "CS263"
$wsB.Range("D5").Value = "MA261"
$wsB.Range("F5").Value = "Free"

$wsB.Range("C6").Value = "Free"

$wsB.Range("B7").Value = "CS261"
$wsB.Range("C7").Value = "CS261"

$wsB.Range("D8").Value = "CS264 (Tutorial)"
$wsB.Range("F8").Value = "CS263 (Tutorial)"
